$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.109894
$ws.Range("H2").Value = 0.329682
$ws.Range("I2").Value = 0.1628895023174853
$ws.Range("J2").Value = 0.1628895023174854
$ws.Range("M2").Value = 0.7435376666666667
$ws.Range("N2").Value = 2.230613
$ws.Range("O2").Value = 0.5505219265933909
$ws.Range("P2").Value = 0.5505219265933909
$ws.Range("Q2").Value = 0.08171032834066666
$ws.Range("R2").Value = 0.7353929550659999
$ws.Range("S2").Value = 0.08967424263766065
$ws.Range("T2").Value = 0.08967424263766066

# Row 3
$ws.Range("G3").Value = 0.109894
$ws.Range("H3").Value = 0.329682
$ws.Range("I3").Value = 0.1628895023174853
$ws.Range("J3").Value = 0.1628895023174854
$ws.Range("O3").Value = 0.08871012126664225
$ws.Range("P3").Value = 0.08871012126664224
$ws.Range("Q3").Value = 0.01316665655933333
$ws.Range("R3").Value = 0.118499909034
$ws.Range("S3").Value = 0.01444994750364713
$ws.Range("T3").Value = 0.01444994750364713

# Row 4
$ws.Range("G4").Value = 0.109894
$ws.Range("H4").Value = 0.329682
$ws.Range("I4").Value = 0.1628895023174853
$ws.Range("J4").Value = 0.1628895023174854
$ws.Range("M4").Value = 0.487255
$ws.Range("N4").Value = 1.461765
$ws.Range("O4").Value = 0.3607679521399669
$ws.Range("P4").Value = 0.3607679521399669
$ws.Range("Q4").Value = 0.05354640097
$ws.Range("R4").Value = 0.4819176087299999
$ws.Range("S4").Value = 0.05876531217617758
$ws.Range("T4").Value = 0.05876531217617759

# Row 5
$ws.Range("I5").Value = 0.6775362766377415
$ws.Range("J5").Value = 0.6775362766377416
$ws.Range("M5").Value = 0.7435376666666667
$ws.Range("N5").Value = 2.230613
$ws.Range("O5").Value = 0.5505219265933909
$ws.Range("P5").Value = 0.5505219265933909
$ws.Range("Q5").Value = 0.3398728023545555
$ws.Range("R5").Value = 3.058855221191
$ws.Range("S5").Value = 0.3729985763515222
$ws.Range("T5").Value = 0.3729985763515222

# Row 6
$ws.Range("I6").Value = 0.6775362766377415
$ws.Range("J6").Value = 0.6775362766377416
$ws.Range("O6").Value = 0.08871012126664225
$ws.Range("P6").Value = 0.08871012126664224
$ws.Range("S6").Value = 0.06010432526308332
$ws.Range("T6").Value = 0.06010432526308333

# Row 7
$ws.Range("I7").Value = 0.6775362766377415
$ws.Range("J7").Value = 0.6775362766377416
$ws.Range("M7").Value = 0.487255
$ws.Range("N7").Value = 1.461765
$ws.Range("O7").Value = 0.3607679521399669
$ws.Range("P7").Value = 0.3607679521399669
$ws.Range("Q7").Value = 0.2227253974283333
$ws.Range("R7").Value = 2.004528576855
$ws.Range("S7").Value = 0.2444333750231361
$ws.Range("T7").Value = 0.2444333750231362

# Row 8
$ws.Range("G8").Value = 0.1076573333333333
$ws.Range("H8").Value = 0.322972
$ws.Range("I8").Value = 0.1595742210447731
$ws.Range("J8").Value = 0.1595742210447731
$ws.Range("M8").Value = 0.7435376666666667
$ws.Range("N8").Value = 2.230613
$ws.Range("O8").Value = 0.5505219265933909
$ws.Range("P8").Value = 0.5505219265933909
$ws.Range("Q8").Value = 0.08004728242622222
$ws.Range("R8").Value = 0.7204255418359999
$ws.Range("S8").Value = 0.08784910760420811
$ws.Range("T8").Value = 0.08784910760420811

# Row 9
$ws.Range("G9").Value = 0.1076573333333333
$ws.Range("H9").Value = 0.322972
$ws.Range("I9").Value = 0.1595742210447731
$ws.Range("J9").Value = 0.1595742210447731
$ws.Range("O9").Value = 0.08871012126664225
$ws.Range("P9").Value = 0.08871012126664224
$ws.Range("Q9").Value = 0.01289867630711111
$ws.Range("R9").Value = 0.116088086764
$ws.Range("S9").Value = 0.0141558484999118
$ws.Range("T9").Value = 0.01415584849991179

# Row 10
$ws.Range("G10").Value = 0.1076573333333333
$ws.Range("H10").Value = 0.322972
$ws.Range("I10").Value = 0.1595742210447731
$ws.Range("J10").Value = 0.1595742210447731
$ws.Range("M10").Value = 0.487255
$ws.Range("N10").Value = 1.461765
$ws.Range("O10").Value = 0.3607679521399669
$ws.Range("P10").Value = 0.3607679521399669
$ws.Range("Q10").Value = 0.05245657395333333
$ws.Range("R10").Value = 0.47210916558
$ws.Range("S10").Value = 0.0575692649406532
$ws.Range("T10").Value = 0.0575692649406532
